$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apos = [string][char]39

$ws.Range('D2').Value = '67.831.07'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '3.807.68'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('D4').Value = $apos + '0.998'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = $apos + '604.76'
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('D6').Value = $apos + '165.94'
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('E10').Value = '  +0.96%  '
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('E12').Value = '  -0.85%  '
$ws.Range('D13').Value = $apos + '36.04'
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').Value = '4.447.46'
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('D15').Value = '3.830.82'
$ws.Range('E15').Value = '  +2.29%  '
$ws.Range('D16').Value = '67.830.23'
$ws.Range('D17').Value = $apos + '18.39'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('E19').Value = '  +1.65%  '
$ws.Range('D20').Value = $apos + '463.76'
$ws.Range('E20').Value = '  +1.30%  '
$ws.Range('D21').Value = $apos + '9.87'
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').Value = $apos + '0.702'
$ws.Range('E22').Value = '  +0.95%  '
$ws.Range('D23').Value = $apos + '0.0000147'
$ws.Range('E23').Value = '  -3.16%  '
$ws.Range('D24').Value = $apos + '83.36'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('E25').Value = '  +0.95%  '
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('D27').Value = $apos + '10.02'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').Value = '3.958.19'
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('E30').Value = '  +0.88%  '
$ws.Range('E31').Value = '  +2.47%  '
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('D33').Value = $apos + '29.46'
$ws.Range('E33').Value = '  -0.62%  '
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').Value = $apos + '9.08'
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('D36').Value = $apos + '0.100'
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').Value = $apos + '0.998'
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('D39').Value = $apos + '5.82'
$ws.Range('E39').Value = '  +1.17%  '
$ws.Range('E40').Value = '  -4.15%  '
$ws.Range('D41').Value = $apos + '0.999'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').Value = $apos + '44.29'
$ws.Range('E43').Value = '  -2.43%  '
$ws.Range('D44').Value = $apos + '47.68'
$ws.Range('E44').Value = '  -0.98%  '
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').Value = $apos + '151.10'
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = $apos + '27.83'
$ws.Range('E47').Value = '  +6.21%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').Value = $apos + '8.36'
$ws.Range('E48').Value = '  +0.51%  '
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').Value = $apos + '1.37'
$ws.Range('E49').Value = '  +11.40%  '
$ws.Range('E50').Value = '  +1.82%  '
$ws.Range('D51').Value = $apos + '388.71'
$ws.Range('E51').Value = '  -0.48%  '
